# Scheduled runner update: refresh cached market-board price/profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4819.727
$ws.Range("J43").Value = 5422.8237
$ws.Range("L43").Value = 5422.8237
$ws.Range("N43").Value = -5560.8237

$ws.Range("H132").Value = 1693.6562
$ws.Range("I132").Value = 1471.8966
$ws.Range("K132").Value = 4415.6898
$ws.Range("M132").Value = -1885.6898

$ws.Range("H137").Value = 3702.652
$ws.Range("I137").Value = 3534.8
$ws.Range("J137").Value = 3749.2778
$ws.Range("K137").Value = 10604.4
$ws.Range("L137").Value = 11247.8334
$ws.Range("M137").Value = -8054.400000000001
$ws.Range("N137").Value = -16347.8334

$ws.Range("H138").Value = 3148.3262
$ws.Range("I138").Value = 2211.0833
$ws.Range("K138").Value = 6633.249899999999
$ws.Range("M138").Value = -1493.249899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 37039250
$ws.Range("I45").Value = 47620276
$ws.Range("K45").Value = 47620276
$ws.Range("M45").Value = -47619899

$ws.Range("H74").Value = 11112649
$ws.Range("I74").Value = 13890452
$ws.Range("K74").Value = 13890452
$ws.Range("M74").Value = -13889578

$ws.Range("H77").Value = 11112649
$ws.Range("I77").Value = 13890452
$ws.Range("K77").Value = 69452260
$ws.Range("M77").Value = -69447892

$ws.Range("H102").Value = 2136.5293
$ws.Range("I102").Value = 2082.5625
$ws.Range("K102").Value = 2082.5625
$ws.Range("M102").Value = -460.5625

$ws.Range("H122").Value = 1432.3
$ws.Range("I122").Value = 659.7857
$ws.Range("K122").Value = 1979.3571
$ws.Range("M122").Value = 470.6428999999998

$ws.Range("H132").Value = 2817.6099
$ws.Range("I132").Value = 2197.2104
$ws.Range("K132").Value = 6591.6312
$ws.Range("M132").Value = -4061.6312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1556.75
$ws.Range("I99").Value = 1180
$ws.Range("K99").Value = 1180
$ws.Range("M99").Value = 318

$ws.Range("H105").Value = 20368.334
$ws.Range("I105").Value = 26679.75
$ws.Range("K105").Value = 26679.75
$ws.Range("M105").Value = -24932.75

$ws.Range("H134").Value = 3015.6843
$ws.Range("I134").Value = 1705.4375
$ws.Range("K134").Value = 5116.3125
$ws.Range("M134").Value = -2581.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2113.8823
$ws.Range("I16").Value = 1390.5
$ws.Range("K16").Value = 1390.5
$ws.Range("M16").Value = -1103.5

$ws.Range("H22").Value = 1771.2142
$ws.Range("I22").Value = 397.14285
$ws.Range("J22").Value = 3145.2856
$ws.Range("K22").Value = 397.14285
$ws.Range("L22").Value = 3145.2856
$ws.Range("M22").Value = -47.14285000000001
$ws.Range("N22").Value = -3845.2856

$ws.Range("H31").Value = 33820.266
$ws.Range("I31").Value = 3055.0952
$ws.Range("J31").Value = 83517.84
$ws.Range("K31").Value = 3055.0952
$ws.Range("L31").Value = 83517.84
$ws.Range("M31").Value = -2760.0952
$ws.Range("N31").Value = -84107.84

$ws.Range("H34").Value = 33820.266
$ws.Range("I34").Value = 3055.0952
$ws.Range("J34").Value = 83517.84
$ws.Range("K34").Value = 3055.0952
$ws.Range("L34").Value = 83517.84
$ws.Range("M34").Value = -2853.0952
$ws.Range("N34").Value = -83921.84

$ws.Range("H107").Value = 905.16
$ws.Range("I107").Value = 811.65
$ws.Range("K107").Value = 811.65
$ws.Range("M107").Value = 1108.35

$ws.Range("H113").Value = 2113.8823
$ws.Range("I113").Value = 1390.5
$ws.Range("K113").Value = 1390.5
$ws.Range("M113").Value = 779.5

$ws.Range("H134").Value = 3943.6428
$ws.Range("I134").Value = 2360.5
$ws.Range("J134").Value = 6054.5
$ws.Range("K134").Value = 7081.5
$ws.Range("L134").Value = 18163.5
$ws.Range("M134").Value = -4546.5
$ws.Range("N134").Value = -23233.5

$ws.Range("H139").Value = 71000
$ws.Range("J139").Value = 71000
$ws.Range("L139").Value = 71000
$ws.Range("N139").Value = -81280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1460.4
$ws.Range("J46").Value = 2750
$ws.Range("L46").Value = 8250
$ws.Range("N46").Value = -8432

$ws.Range("H113").Value = 1515.4166
$ws.Range("I113").Value = 663.6667
$ws.Range("J113").Value = 1799.3334
$ws.Range("K113").Value = 1991.0001
$ws.Range("L113").Value = 5398.0002
$ws.Range("M113").Value = 178.9999
$ws.Range("N113").Value = -9738.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 19346.5
$ws.Range("J33").Value = 19346.5
$ws.Range("L33").Value = 19346.5
$ws.Range("N33").Value = -19850.5

$ws.Range("H97").Value = 1167.8636
$ws.Range("I97").Value = 665.9167
$ws.Range("J97").Value = 1770.2
$ws.Range("K97").Value = 665.9167
$ws.Range("L97").Value = 1770.2
$ws.Range("M97").Value = -169.9167
$ws.Range("N97").Value = -2762.2

$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H113").Value = 3711.2307
$ws.Range("I113").Value = 3023.2727
$ws.Range("K113").Value = 3023.2727
$ws.Range("M113").Value = -853.2727

$ws.Range("H132").Value = 39259.535
$ws.Range("I132").Value = 41086.652
$ws.Range("K132").Value = 123259.956
$ws.Range("M132").Value = -120729.956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 23512
$ws.Range("J34").Value = 23512
$ws.Range("L34").Value = 23512
$ws.Range("N34").Value = -23856

$ws.Range("H64").Value = 39000
$ws.Range("J64").Value = 39000
$ws.Range("L64").Value = 39000
$ws.Range("N64").Value = -39450

$ws.Range("H67").Value = 39000
$ws.Range("J67").Value = 39000
$ws.Range("L67").Value = 39000
$ws.Range("N67").Value = -40560

$ws.Range("H74").Value = 55000
$ws.Range("J74").Value = 55000
$ws.Range("L74").Value = 55000
$ws.Range("N74").Value = -56996

$ws.Range("H77").Value = 55000
$ws.Range("J77").Value = 55000
$ws.Range("L77").Value = 165000
$ws.Range("N77").Value = -174984

$ws.Range("H132").Value = 4782.4062
$ws.Range("I132").Value = 2521.32
$ws.Range("K132").Value = 7563.960000000001
$ws.Range("M132").Value = -5033.960000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 1250
$ws.Range("J31").Value = 1250
$ws.Range("L31").Value = 1250
$ws.Range("N31").Value = -1946

$ws.Range("H62").Value = 9158
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 9158
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240

$ws.Range("H96").Value = 3783.3333
$ws.Range("I96").Value = 2966.8333
$ws.Range("K96").Value = 2966.8333
$ws.Range("M96").Value = -1593.8333

$ws.Range("H131").Value = 76666.664
$ws.Range("J131").Value = 76666.664
$ws.Range("L131").Value = 76666.664
$ws.Range("N131").Value = -86746.664

